# Apply "I0 and IF added" edit: add two new columns (I = "I0", J = "IF")
# with per-row numeric data, mirroring the style of the existing header
# cell H1 and extending the sheet's used range to A1:J75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers I1 = "I0", J1 = "IF" ---
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Copy the header style (bold, bordered, centered) used by the existing
# header cells (e.g. H1) onto the two new header cells.
$ws.Cells.Item(1, 8).Copy() | Out-Null
$ws.Cells.Item(1, 9).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(1, 10).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Data rows 2-75: numeric values for I0 and IF columns ---
$i0Values = @(8,8,8,9,8,8,8,7,9,8,8,7,3,8,7,8,8,8,8,9,7,8,8,7,8,7,8,8,8,7,7,8,8,8,9,8,8,8,8,8,8,8,8,8,11,8,9,8,8,8,5,7,8,8,8,6,9,7,7,9,7,8,8,8,7,8,7,6,5,4,4,4,6,4)
$ifValues = @(8,8,8,9,8,8,8,7,9,8,8,8,4,8,7,8,8,8,8,9,7,8,8,7,8,8,8,8,8,7,8,8,8,8,9,8,8,8,8,8,8,8,8,8,11,8,9,8,8,8,5,7,8,8,8,6,9,7,7,9,7,8,8,8,7,8,7,7,5,4,4,4,6,4)

for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}
